# Re-pull / push data: update the dSF column (F) values for the affected rows.
# These cells previously mirrored the dS0 column (E) but after a repull/mean
# recalculation they now hold distinct values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -6
    5  = -7
    7  = 1
    10 = -2
    13 = -10
    17 = -9
    18 = 4
    20 = -5
    21 = -5
    22 = -4
    27 = -1
    32 = -9
    33 = -1
    34 = -3
    37 = -3
    50 = 0
    51 = 3
    56 = 2
    57 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
